# Adds 4 new data rows (19-22) to the "Artfynd" sheet, extending the
# used range from A1:AY18 to A1:AY22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range('A19').Value = 111926769
$ws.Range('B19').Value = 90658
$ws.Range('C19').Value = 'Ovaliderad'
$ws.Range('D19').Value = 'NT'
$ws.Range('E19').Value = 4361
$ws.Range('F19').Value = 'Orange taggsvamp'
$ws.Range('G19').Value = 'Hydnellum aurantiacum'
$ws.Range('H19').Value = '(Batsch:Fr.) P.Karst.'
$ws.Range('I19').Value = "'"
$ws.Range('K19').Value = "'"
$ws.Range('P19').Value = 'Upplands Väsby (Upplands Väsby), Upl'
$ws.Range('Q19').Value = 663476.3239106013
$ws.Range('R19').Value = 6602651.048317727
$ws.Range('S19').Value = 10
$ws.Range('T19').Value = 'Stockholm'
$ws.Range('U19').Value = 'Upplands Väsby'
$ws.Range('V19').Value = 'Uppland'
$ws.Range('W19').Value = 'Ed'
$ws.Range('Y19').Value = "'2023-09-06"
$ws.Range('Z19').Value = '00:00'
$ws.Range('AA19').Value = "'2023-09-06"
$ws.Range('AB19').Value = '00:00'
$ws.Range('AD19').Value = $false
$ws.Range('AE19').Value = $false
$ws.Range('AG19').Value = $false
$ws.Range('AT19').Value = "'"
$ws.Range('AW19').Value = 'Hans Bärring'
$ws.Range('AX19').Value = 'Hans Bärring'
$ws.Range('AY19').Value = "'"

# Row 20
$ws.Range('A20').Value = 111929648
$ws.Range('B20').Value = 90655
$ws.Range('C20').Value = 'Ovaliderad'
$ws.Range('D20').Value = 'VU'
$ws.Range('E20').Value = 150
$ws.Range('F20').Value = 'Grangråticka'
$ws.Range('G20').Value = 'Boletopsis leucomelaena'
$ws.Range('H20').Value = '(Pers.) Fayod'
$ws.Range('I20').Value = "'10"
$ws.Range('J20').Value = 'fruktkroppar'
$ws.Range('K20').Value = "'"
$ws.Range('P20').Value = 'Upplands Väsby (Upplands Väsby), Upl'
$ws.Range('Q20').Value = 663509.4011657666
$ws.Range('R20').Value = 6602732.737488487
$ws.Range('S20').Value = 10
$ws.Range('T20').Value = 'Stockholm'
$ws.Range('U20').Value = 'Upplands Väsby'
$ws.Range('V20').Value = 'Uppland'
$ws.Range('W20').Value = 'Ed'
$ws.Range('Y20').Value = "'2023-09-06"
$ws.Range('Z20').Value = '00:00'
$ws.Range('AA20').Value = "'2023-09-06"
$ws.Range('AB20').Value = '00:00'
$ws.Range('AD20').Value = $false
$ws.Range('AE20').Value = $false
$ws.Range('AG20').Value = $false
$ws.Range('AT20').Value = "'"
$ws.Range('AW20').Value = 'Hans Bärring'
$ws.Range('AX20').Value = 'Hans Bärring'
$ws.Range('AY20').Value = "'"

# Row 21
$ws.Range('A21').Value = 111927215
$ws.Range('B21').Value = 90658
$ws.Range('C21').Value = 'Ovaliderad'
$ws.Range('D21').Value = 'NT'
$ws.Range('E21').Value = 4361
$ws.Range('F21').Value = 'Orange taggsvamp'
$ws.Range('G21').Value = 'Hydnellum aurantiacum'
$ws.Range('H21').Value = '(Batsch:Fr.) P.Karst.'
$ws.Range('I21').Value = "'10"
$ws.Range('J21').Value = 'fruktkroppar'
$ws.Range('K21').Value = "'"
$ws.Range('P21').Value = 'Upplands Väsby (Upplands Väsby), Upl'
$ws.Range('Q21').Value = 663485.6413922446
$ws.Range('R21').Value = 6602647.390513759
$ws.Range('S21').Value = 10
$ws.Range('T21').Value = 'Stockholm'
$ws.Range('U21').Value = 'Upplands Väsby'
$ws.Range('V21').Value = 'Uppland'
$ws.Range('W21').Value = 'Ed'
$ws.Range('Y21').Value = "'2023-09-06"
$ws.Range('Z21').Value = '00:00'
$ws.Range('AA21').Value = "'2023-09-06"
$ws.Range('AB21').Value = '00:00'
$ws.Range('AD21').Value = $false
$ws.Range('AE21').Value = $false
$ws.Range('AG21').Value = $false
$ws.Range('AT21').Value = "'"
$ws.Range('AW21').Value = 'Hans Bärring'
$ws.Range('AX21').Value = 'Hans Bärring'
$ws.Range('AY21').Value = "'"

# Row 22
$ws.Range('A22').Value = 111926622
$ws.Range('B22').Value = 90658
$ws.Range('C22').Value = 'Ovaliderad'
$ws.Range('D22').Value = 'NT'
$ws.Range('E22').Value = 4361
$ws.Range('F22').Value = 'Orange taggsvamp'
$ws.Range('G22').Value = 'Hydnellum aurantiacum'
$ws.Range('H22').Value = '(Batsch:Fr.) P.Karst.'
$ws.Range('I22').Value = "'"
$ws.Range('K22').Value = "'"
$ws.Range('P22').Value = 'Upplands Väsby (Upplands Väsby), Upl'
$ws.Range('Q22').Value = 663452.3464515609
$ws.Range('R22').Value = 6602675.90838708
$ws.Range('S22').Value = 10
$ws.Range('T22').Value = 'Stockholm'
$ws.Range('U22').Value = 'Upplands Väsby'
$ws.Range('V22').Value = 'Uppland'
$ws.Range('W22').Value = 'Ed'
$ws.Range('Y22').Value = "'2023-09-06"
$ws.Range('Z22').Value = '00:00'
$ws.Range('AA22').Value = "'2023-09-06"
$ws.Range('AB22').Value = '00:00'
$ws.Range('AD22').Value = $false
$ws.Range('AE22').Value = $false
$ws.Range('AG22').Value = $false
$ws.Range('AT22').Value = "'"
$ws.Range('AW22').Value = 'Hans Bärring'
$ws.Range('AX22').Value = 'Hans Bärring'
$ws.Range('AY22').Value = "'"
